$d = $word.ActiveDocument
$d.Content.Find.Execute("Hello this is version 1.", $true, $false, $false, $false, $false, $true, 1, $false, "Hello this is version 2.", 2)
